$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.853403627872467
$ws.Range("B1").Value = 1.216473579406738
$ws.Range("C1").Value = 3.41157603263855
$ws.Range("D1").Value = 4.093273162841797
$ws.Range("E1").Value = 0.6605732440948486
